$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new students appended to the roster (rows 22 and 23).
$ws.Range("A22").Value2 = "Đặng Ngọc Anh"
$ws.Range("B22").Value2 = 215748020110333
$ws.Range("C22").Value2 = 344483272

$ws.Range("A23").Value2 = "Thái Văn Tuấn"
$ws.Range("B23").Value2 = 215748020110086
$ws.Range("C23").Value2 = 987654321

# Match the author's final selection/scroll position in the sheet.
$ws.Range("E23").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
